# Update "想去人数" (want-to-go count) figures in column F for both the
# "展览" sheet and the "全部类型" sheet, per the latest scrape refresh.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - rows keyed by event, F column updates
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 14926
$wsExpo.Range("F3").Value = 18675
$wsExpo.Range("F5").Value = 121
$wsExpo.Range("F15").Value = 205
$wsExpo.Range("F22").Value = 7771

# Sheet "全部类型" (all types) - same events, slightly shifted row numbers
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 14926
$wsAll.Range("F3").Value = 18675
$wsAll.Range("F5").Value = 121
$wsAll.Range("F15").Value = 205
$wsAll.Range("F23").Value = 7771
